# Reassign the "usuario" (user) value for rows 26-31 in the MITIGACION sheet
# from "matias.larenti" to the new value "eduardo.XXXXXXXX", and update the
# sheet selection to reflect the edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MITIGACION")

$ws.Range("A26:A31").Value = "eduardo.XXXXXXXX"

$ws.Range("A26:A31").Select()
